$d = $word.ActiveDocument

# ---- "Programa" paragraph ----
# Original body text is one big run; split it into separate sentences
# joined by manual line breaks (^l => <w:br/>), matching the target
# formatting (one w:r containing multiple w:t/w:br children).
$programaOld = "Conteúdo teórico: 1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. 5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.Conteúdo prático: 1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.2. Visita a fabricantes de abrasivos. 3. Visita a usuários de ferramentas abrasivas."

$programaNew = "Conteúdo teórico: ^l1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.^l2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.^l3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.^l4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. ^l5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.^l6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.^l^lConteúdo prático: ^l1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.^l2. Visita a fabricantes de abrasivos. ^l3. Visita a usuários de ferramentas abrasivas."

$range1 = $d.Content
$found1 = $range1.Find.Execute($programaOld, $true, $true, $false, $false, $false, $true, 1, $false, $programaNew, 2)
if (-not $found1) {
    throw "Could not find/replace the 'Programa' paragraph text"
}

# ---- "Bibliografia" paragraph ----
# Same treatment: split the single run of concatenated references into
# one reference per line via manual line breaks.
$biblioOld = "1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. 2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004."

$biblioNew = "1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. ^l2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.^l3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.^l4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.^l5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.^l6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.^l7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.^l8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.^l9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004."

$range2 = $d.Content
$found2 = $range2.Find.Execute($biblioOld, $true, $true, $false, $false, $false, $true, 1, $false, $biblioNew, 2)
if (-not $found2) {
    throw "Could not find/replace the 'Bibliografia' paragraph text"
}

Write-Output "Programa replaced: $found1; Bibliografia replaced: $found2"
